# Generate Report for Handoff
# "b.md" has now been handed off: status flips from
# "Handed back: in sync with en-US" to "Ready for handoff" everywhere it is
# reported, a new (later) handback xliff is recorded, and an error is
# surfaced because that handback isn't the latest version.

$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/oltest/blob/ea390eaa40dc00dc8afc058ae26070a94ecf6d1c/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/oltest/blob/df4b611567d87ab84d073ba2e8f2fc69d0f13c10/e2e/b.md."

# ---- Overview sheet: row 3 is the b.md row ----
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = "Ready for handoff"
$overview.Range("F3").Value = "Ready for handoff"
$overview.Range("G3").Value = "2016-08-05 02:34:47"

# ---- zh-cn sheet: row 3 is the b.md row ----
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Columns.Item(16).ColumnWidth = 39.17
$zhcn.Range("C3").Value = "Ready for handoff"
# Copy the literal text "False" from another cell on the same row so Excel
# keeps it as text instead of auto-converting the word into a boolean.
$zhcn.Range("O3").Copy($zhcn.Range("F3"))
$zhcn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zhcn.Range("H3").Value = "2016-08-05 02:34:37"
$zhcn.Range("P3").Value = $errorDetail

# ---- de-de sheet: row 3 is the b.md row ----
$dede = $wb.Worksheets.Item("de-de")
$dede.Columns.Item(16).ColumnWidth = 39.17
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("O3").Copy($dede.Range("F3"))
$dede.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$dede.Range("H3").Value = "2016-08-05 02:34:47"
$dede.Range("P3").Value = $errorDetail
